# Rename the sheet "Unidade DAS_CTC" -> "DasCTC"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DasCTC"

# Move the active cell selection to H9 (column "Aparelho gela?")
$ws.Range("H9").Select()
